# Corrected the global models table: round/format p-values and CIs.
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "1.873034e-05"; New = "< 0.001" },
    @{ Old = "1.081093e-03"; New = "0.001" },
    @{ Old = "[0.217, 0.695]"; New = "[0.22, 0.70]" },
    @{ Old = "6.363"; New = "6.36" },
    @{ Old = "5.291"; New = "5.29" },
    @{ Old = "1.536603e-03"; New = "0.002" },
    @{ Old = "[0.066, 0.274]"; New = "[0.07, 0.27]" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
